$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

for ($i = 1; $i -le 8; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = "Announcement$i-230120"
    $ws.Cells.Item($row, 2).Value = "Announcement$i-230120 Summary"
}
